$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update values
$ws.Range("B3").Value = "Unknown"
$ws.Range("B4").Value = "22-06-2020"
$ws.Range("B5").Value = 3

# Row 10: Type column (D10) and Attribute description (E10)
$ws.Range("D10").Value = "String"
$descText = "Values: `nBird = Aves.`nBat = Microchiroptera. "
$ws.Range("E10").Value = $descText

$run1 = $ws.Range("E10").Characters(1, 7)
$run1.Font.Underline = $true

$run2Len = $descText.Length - 7
$run2 = $ws.Range("E10").Characters(8, $run2Len)
$run2.Font.Name = "Times New Roman"
$run2.Font.Size = 9

# Row 10 grew taller to fit the extra wrapped text
$ws.Range("A10").EntireRow.RowHeight = 36.75

# Page layout tweak
$ws.PageSetup.Orientation = 1

$ws.Range("B4").Select()
